# Update the "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets. Each sheet receives the same set of +1 (and one +2) bumps.

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 621
    10 = 398
    12 = 117
    19 = 315
    22 = 88
    26 = 244
    29 = 1658
    35 = 3864
    40 = 82
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
